{"js": "// Word JS API (Office.js) edit script.\n// Body is the async (context) => { ... } function content.\n//\n// The authored change (AQL 7.x -> 8.x) updates the wording of an\n// \"invalid expression\" diagnostic that M2Doc stamps into the document\n// when a `{m:link ...}` field resolves to an empty expression:\n//\n//   Expression \"\" is invalid: null or empty string.\n//     -> Expression \"\" is invalid: missing expression\n//\n// Find the run(s) containing the old diagnostic text and replace the\n// text in place so surrounding run formatting (bold + red colour) is\n// preserved.\nconst oldText = 'Expression \"\" is invalid: null or empty string.';\nconst newText = 'Expression \"\" is invalid: missing expression';\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# $word / $d (ActiveDocument) are pre-seeded by the harness.\n#\n# The authored change (AQL 7.x -> 8.x) updates the wording of an\n# \"invalid expression\" diagnostic that M2Doc stamps into the document\n# when a `{m:link ...}` field resolves to an empty expression:\n#\n#   Expression \"\" is invalid: null or empty string.\n#     -> Expression \"\" is invalid: missing expression\n#\n# Use Find/Replace on the document content so the surrounding run\n# formatting (bold + red colour) is preserved.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = \"is invalid: null or empty string.\"\n$find.Replacement.Text = \"is invalid: missing expression\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
